$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.724031686782837
$ws.Range("B1").Value = 2.392427444458008
$ws.Range("C1").Value = 2.644598245620728
$ws.Range("D1").Value = 3.306588888168335
$ws.Range("E1").Value = 2.942408084869385
